$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Manchester tribunal details
$ws.Range("B3").Value  = "Manchester Employment Tribunal, Alexandra House, 14-22 The Parsonage, Manchester, M3 2JA"
$ws.Range("B4").Value  = "0161 833 6100"
$ws.Range("B5").Value  = "0870 739 4433"
$ws.Range("B6").Value  = "DX 743570"
$ws.Range("B7").Value  = "Manchesteret@justice.gov.uk"

# Glasgow tribunal details
$ws.Range("B8").Value  = "Eagle Building, 215 Bothwell Street, Glasgow, G2 7TS"
$ws.Range("B9").Value  = "0141 204 0730"
$ws.Range("B10").Value = "01264 785 177"
$ws.Range("B11").Value = "DX 7435701"
$ws.Range("B12").Value = "glasgowet@justice.gov.uk"

# Keep the mailto hyperlinks' visible text in sync with the new cell values
foreach ($h in $ws.Hyperlinks) {
    if ($h.Address -eq "mailto:manchester@gmail.com") {
        $h.TextToDisplay = "Manchesteret@justice.gov.uk"
    }
    elseif ($h.Address -eq "mailto:glasgow@gmail.com") {
        $h.TextToDisplay = "glasgowet@justice.gov.uk"
    }
}

Write-Output "venue details updated"
